$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.845.13"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.713.00"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'318.19"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "'0.3970"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").Value = "'1.530"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "'1.005"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").Value = "'53.44"
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").Value = "'0.08967"
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "'7.697"
$ws.Range("E13").Value = "  +6.54%  "
$ws.Range("D14").Value = "'24.40"
$ws.Range("E14").Value = "  +4.50%  "
$ws.Range("D15").Value = "'8.198"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "'0.00001363"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "1.705.42"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "'100.35"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'0.07160"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "'20.28"
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Value = "'7.514"
$ws.Range("E21").Value = "  +6.12%  "
$ws.Range("D22").Value = "'1.006"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "'14.58"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").Value = "24.837.37"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "'3.112"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").Value = "'2.344"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'23.09"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "'9.340"
$ws.Range("E28").Value = "  +23.88%  "
$ws.Range("D29").Value = "'166.52"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "'139.72"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "'5.239"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "'7.925"
$ws.Range("E32").Value = "  +10.11%  "
$ws.Range("D33").Value = "'0.09105"
$ws.Range("B34").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C34").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D34").Value = "1.892.65"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.097"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.03027"
$ws.Range("E36").Value = "  +10.67%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2820"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'11.15"
$ws.Range("E38").Value = "  -4.02%  "
$ws.Range("B39").Value = "WEMIXTOKEN"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.963"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'0.09353"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.8168"
$ws.Range("E42").Value = "  +6.64%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.491"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'16.89"
$ws.Range("E44").Value = "  +8.15%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.7432"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.661"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'4.278"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "'1.360"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "'1.002"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'141.33"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'94.46"
$ws.Range("E51").Value = "  +5.20%  "
